# Update the "想去人数" (wanted-to-go count) column F values on the
# "展览" (Exhibition) and "全部类型" (All types) sheets to reflect the
# refreshed gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 543
$ws1.Range("F4").Value = 1562
$ws1.Range("F5").Value = 163
$ws1.Range("F8").Value = 178
$ws1.Range("F9").Value = 758
$ws1.Range("F10").Value = 1053
$ws1.Range("F12").Value = 349
$ws1.Range("F13").Value = 62
$ws1.Range("F14").Value = 102
$ws1.Range("F15").Value = 18
$ws1.Range("F16").Value = 6510
$ws1.Range("F17").Value = 25
$ws1.Range("F20").Value = 161
$ws1.Range("F22").Value = 15544
$ws1.Range("F23").Value = 1537
$ws1.Range("F24").Value = 293
$ws1.Range("F25").Value = 150
$ws1.Range("F27").Value = 11110
$ws1.Range("F28").Value = 773
$ws1.Range("F29").Value = 4355
$ws1.Range("F33").Value = 309
$ws1.Range("F34").Value = 129

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 543
$ws4.Range("F4").Value = 1562
$ws4.Range("F5").Value = 163
$ws4.Range("F9").Value = 178
$ws4.Range("F10").Value = 758
$ws4.Range("F12").Value = 1053
$ws4.Range("F14").Value = 349
$ws4.Range("F15").Value = 62
$ws4.Range("F16").Value = 104
$ws4.Range("F18").Value = 18
$ws4.Range("F19").Value = 6510
$ws4.Range("F20").Value = 25
$ws4.Range("F23").Value = 161
$ws4.Range("F26").Value = 15544
$ws4.Range("F27").Value = 1537
$ws4.Range("F28").Value = 293
$ws4.Range("F29").Value = 150
$ws4.Range("F32").Value = 11110
$ws4.Range("F33").Value = 774
$ws4.Range("F34").Value = 4355
$ws4.Range("F37").Value = 21
$ws4.Range("F38").Value = 309
$ws4.Range("F39").Value = 129
